$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '22.427.07'
$r.Style = 'Normal'
$r = $ws.Range('D3')
$r.NumberFormat = '@'
$r.Value = '1.571.47'
$r.Style = 'Normal'
$r = $ws.Range('E4')
$r.NumberFormat = '@'
$r.Value = '  +0.08%  '
$r.Style = 'Normal'
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '1.003'
$r.Style = 'Normal'
$r = $ws.Range('E5')
$r.NumberFormat = '@'
$r.Value = '  +0.06%  '
$r.Style = 'Normal'
$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '291.36'
$r.Style = 'Normal'
$r = $ws.Range('E6')
$r.NumberFormat = '@'
$r.Value = '  +0.31%  '
$r.Style = 'Normal'
$r = $ws.Range('D7')
$r.NumberFormat = '@'
$r.Value = '0.3727'
$r.Style = 'Normal'
$r = $ws.Range('E7')
$r.NumberFormat = '@'
$r.Value = '  -1.05%  '
$r.Style = 'Normal'
$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '50.08'
$r.Style = 'Normal'
$r = $ws.Range('E8')
$r.NumberFormat = '@'
$r.Value = '  +0.37%  '
$r.Style = 'Normal'
$r = $ws.Range('D9')
$r.NumberFormat = '@'
$r.Value = '0.3387'
$r.Style = 'Normal'
$r = $ws.Range('E9')
$r.NumberFormat = '@'
$r.Value = '  -1.03%  '
$r.Style = 'Normal'
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '0.07576'
$r.Style = 'Normal'
$r = $ws.Range('E10')
$r.NumberFormat = '@'
$r.Value = '  -0.95%  '
$r.Style = 'Normal'
$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '1.138'
$r.Style = 'Normal'
$r = $ws.Range('E11')
$r.NumberFormat = '@'
$r.Value = '  -1.66%  '
$r.Style = 'Normal'
$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '1.002'
$r.Style = 'Normal'
$r = $ws.Range('E12')
$r.NumberFormat = '@'
$r.Value = '  -0.07%  '
$r.Style = 'Normal'
$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '21.31'
$r.Style = 'Normal'
$r = $ws.Range('E13')
$r.NumberFormat = '@'
$r.Value = '  +0.25%  '
$r.Style = 'Normal'
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '5.986'
$r.Style = 'Normal'
$r = $ws.Range('E14')
$r.NumberFormat = '@'
$r.Value = '  -0.69%  '
$r.Style = 'Normal'
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '6.956'
$r.Style = 'Normal'
$r = $ws.Range('E15')
$r.NumberFormat = '@'
$r.Value = '  +0.10%  '
$r.Style = 'Normal'
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '1.573.24'
$r.Style = 'Normal'
$r = $ws.Range('E16')
$r.NumberFormat = '@'
$r.Value = '  +0.04%  '
$r.Style = 'Normal'
$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '0.00001119'
$r.Style = 'Normal'
$r = $ws.Range('E17')
$r.NumberFormat = '@'
$r.Value = '  -1.27%  '
$r.Style = 'Normal'
$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '90.92'
$r.Style = 'Normal'
$r = $ws.Range('E18')
$r.NumberFormat = '@'
$r.Value = '  +0.76%  '
$r.Style = 'Normal'
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '0.06733'
$r.Style = 'Normal'
$r = $ws.Range('E19')
$r.NumberFormat = '@'
$r.Value = '  -0.40%  '
$r.Style = 'Normal'
$r = $ws.Range('E20')
$r.NumberFormat = '@'
$r.Value = '  +0.15%  '
$r.Style = 'Normal'
$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '6.284'
$r.Style = 'Normal'
$r = $ws.Range('E21')
$r.NumberFormat = '@'
$r.Value = '  +1.23%  '
$r.Style = 'Normal'
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '16.31'
$r.Style = 'Normal'
$r = $ws.Range('E22')
$r.NumberFormat = '@'
$r.Value = '  -3.08%  '
$r.Style = 'Normal'
$r = $ws.Range('E23')
$r.NumberFormat = '@'
$r.Value = '  +1.05%  '
$r.Style = 'Normal'
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '22.448.43'
$r.Style = 'Normal'
$r = $ws.Range('E24')
$r.NumberFormat = '@'
$r.Value = '  +0.24%  '
$r.Style = 'Normal'
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '2.336'
$r.Style = 'Normal'
$r = $ws.Range('E25')
$r.NumberFormat = '@'
$r.Value = '  -3.78%  '
$r.Style = 'Normal'
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '2.662'
$r.Style = 'Normal'
$r = $ws.Range('E26')
$r.NumberFormat = '@'
$r.Value = '  -1.38%  '
$r.Style = 'Normal'
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '20.12'
$r.Style = 'Normal'
$r = $ws.Range('E27')
$r.NumberFormat = '@'
$r.Value = '  -0.70%  '
$r.Style = 'Normal'
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '148.43'
$r.Style = 'Normal'
$r = $ws.Range('E28')
$r.NumberFormat = '@'
$r.Value = '  +0.76%  '
$r.Style = 'Normal'
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '5.013'
$r.Style = 'Normal'
$r = $ws.Range('E29')
$r.NumberFormat = '@'
$r.Value = '  -0.53%  '
$r.Style = 'Normal'
$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '125.39'
$r.Style = 'Normal'
$r = $ws.Range('E30')
$r.NumberFormat = '@'
$r.Value = '  -0.80%  '
$r.Style = 'Normal'
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '1.750.56'
$r.Style = 'Normal'
$r = $ws.Range('E31')
$r.NumberFormat = '@'
$r.Value = '  +0.14%  '
$r.Style = 'Normal'
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '1.046'
$r.Style = 'Normal'
$r = $ws.Range('E32')
$r.NumberFormat = '@'
$r.Value = '  +4.98%  '
$r.Style = 'Normal'
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '6.158'
$r.Style = 'Normal'
$r = $ws.Range('E33')
$r.NumberFormat = '@'
$r.Value = '  -0.14%  '
$r.Style = 'Normal'
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '1.973'
$r.Style = 'Normal'
$r = $ws.Range('E34')
$r.NumberFormat = '@'
$r.Value = '  -1.74%  '
$r.Style = 'Normal'
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '9.808'
$r.Style = 'Normal'
$r = $ws.Range('E35')
$r.NumberFormat = '@'
$r.Value = '  -1.78%  '
$r.Style = 'Normal'
$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '0.08379'
$r.Style = 'Normal'
$r = $ws.Range('E36')
$r.NumberFormat = '@'
$r.Value = '  -2.30%  '
$r.Style = 'Normal'
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '1.377'
$r.Style = 'Normal'
$r = $ws.Range('E37')
$r.NumberFormat = '@'
$r.Value = '  +3.83%  '
$r.Style = 'Normal'
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '0.02472'
$r.Style = 'Normal'
$r = $ws.Range('E38')
$r.NumberFormat = '@'
$r.Value = '  -3.05%  '
$r.Style = 'Normal'
$r = $ws.Range('E39')
$r.NumberFormat = '@'
$r.Value = '  -1.53%  '
$r.Style = 'Normal'
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '0.06521'
$r.Style = 'Normal'
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '5.454'
$r.Style = 'Normal'
$r = $ws.Range('E41')
$r.NumberFormat = '@'
$r.Value = '  +0.03%  '
$r.Style = 'Normal'
$r = $ws.Range('E42')
$r.NumberFormat = '@'
$r.Value = '  -2.31%  '
$r.Style = 'Normal'
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '0.6223'
$r.Style = 'Normal'
$r = $ws.Range('E43')
$r.NumberFormat = '@'
$r.Value = '  -3.11%  '
$r.Style = 'Normal'
$r = $ws.Range('E44')
$r.NumberFormat = '@'
$r.Value = '  +0.07%  '
$r.Style = 'Normal'
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '14.00'
$r.Style = 'Normal'
$r = $ws.Range('E45')
$r.NumberFormat = '@'
$r.Value = '  +0.00%  '
$r.Style = 'Normal'
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '3.812'
$r.Style = 'Normal'
$r = $ws.Range('E46')
$r.NumberFormat = '@'
$r.Value = '  +0.52%  '
$r.Style = 'Normal'
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '0.5787'
$r.Style = 'Normal'
$r = $ws.Range('E47')
$r.NumberFormat = '@'
$r.Value = '  -3.52%  '
$r.Style = 'Normal'
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '129.42'
$r.Style = 'Normal'
$r = $ws.Range('E48')
$r.NumberFormat = '@'
$r.Value = '  +2.99%  '
$r.Style = 'Normal'
$r = $ws.Range('D49')
$r.NumberFormat = '@'
$r.Value = '2.072'
$r.Style = 'Normal'
$r = $ws.Range('E49')
$r.NumberFormat = '@'
$r.Value = '  -0.73%  '
$r.Style = 'Normal'
$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '1.213'
$r.Style = 'Normal'
$r = $ws.Range('E50')
$r.NumberFormat = '@'
$r.Value = '  -6.91%  '
$r.Style = 'Normal'
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '0.07317'
$r.Style = 'Normal'
$r = $ws.Range('E51')
$r.NumberFormat = '@'
$r.Value = '  -0.19%  '
$r.Style = 'Normal'
